# Leave Card update — adds a new "SL(1-0-0)" leave-card row (table row 95 / sheet
# row 102), pushing the existing trailing template rows down by one, and records
# an "A(1-0-0)" remark date "8/23,24/2023" against the VL(2-0-0) row above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item("Table1")

# 1) Duplicate the current bottom "totals" style row (147, special border style)
#    down onto the new physical last row (148), so the table keeps its closing
#    bottom border after it grows by one row.
$ws.Range("A147:K147").Copy($ws.Range("A148:K148"))

# 2) The row that used to be the special bottom-border row (147) becomes a
#    normal blank body row, matching the row above it (146).
$ws.Range("A146:K146").Copy($ws.Range("A147:K147"))

# 3) Shift the PERIOD (date) column down by one row for every blank template
#    row between the insertion point and the old last row, so each row keeps
#    the date that used to belong to the row above it.
$ws.Range("A102:A141").Copy($ws.Range("A103:A142"))

# 4) Populate the newly inserted row (now sheet row 102 / table row 95) with
#    the new SL(1-0-0) leave entry.
$ws.Range("A102").Value = $null
$ws.Range("B102").Value = "SL(1-0-0)"
$ws.Range("C102").Value = $null
$ws.Range("D102").Value = $null
$ws.Range("H102").Value = 1
$ws.Range("K102").Style = $ws.Range("K90").Style
$ws.Range("K102").Value = 45152

# 5) Fill in the EARNED figures that were recorded for the surrounding rows.
$ws.Range("C98").Value = 1.25
$ws.Range("C100").Value = 1.25

# 6) Record the new VL(2-0-0) leave (2 days) and its "A(1-0-0)" remark date.
$ws.Range("B101").Value = "VL(2-0-0)"
$ws.Range("D101").Value = 2
$ws.Range("K101").Value = "8/23,24/2023"

# 7) Grow the table to include the new trailing row.
$tbl.Resize($ws.Range("A8:K148"))

# 8) Restore the view state (active cell / scroll position) to match.
$ws.Range("K102").Select()
$excel.ActiveWindow.ScrollRow = 91
